# Added PSA talk recording.
# Inserts a new "Meetups" row for the 2023-04-24 PSA (Estimating Causality
# from Observational Data) talk, complete with its recording link and
# slides/notes, pushing the existing later meetups down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Meetups")

# Insert a fresh row at 15 (copies formatting down from row 14, shifts the
# old rows 15-18 down to 16-19).
$ws.Rows.Item(15).Insert()

$ws.Cells.Item(15, 1).Value = 45040
$ws.Cells.Item(15, 2).Value = "6:30 pm"
$ws.Cells.Item(15, 3).Value = "8:00 pm"
$ws.Cells.Item(15, 4).Value = "Estimating Causality from Observational Data"
$ws.Cells.Item(15, 5).Value = "/course-overview/meetups"
$ws.Cells.Item(15, 6).Value = "Rq_od5KwqEA"
$ws.Cells.Item(15, 8).Value = "Slides: https://github.com/jbryer/psa/raw/master/Slides/Intro_PSA.pdf <br/>`nBookdown site: https://psa.bryer.org <br/>`nGithub repo: https://github.com/jbryer/psa"
$ws.Cells.Item(15, 8).WrapText = $true
$ws.Cells.Item(15, 8).Font.Color = 0
$ws.Rows.Item(15).RowHeight = 85

# Match the saved selection from the source edit.
$ws.Range("H16").Select()

Write-Output "PSA talk recording row inserted."
